# Apply updated cryptocurrency price/volume data to the worksheet.
# Values are set with a leading apostrophe so that Excel stores them as
# literal text (matching the original inlineStr/text cells) instead of
# auto-converting number-like strings (e.g. "26.110.36", "0.06349") into
# numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.110.36"
$ws.Range("E2").Value = "'  -0.23%  "
# Row 3
$ws.Range("D3").Value = "'1.655.20"
$ws.Range("E3").Value = "'  -0.17%  "
# Row 4
$ws.Range("E4").Value = "'  -0.34%  "
# Row 5
$ws.Range("E5").Value = "'  +0.45%  "
# Row 6
$ws.Range("E6").Value = "'  +1.18%  "
# Row 7
$ws.Range("E7").Value = "'  -0.28%  "
# Row 8
$ws.Range("E8").Value = "'  -0.79%  "
# Row 9
$ws.Range("D9").Value = "'0.06349"
$ws.Range("E9").Value = "'  +1.88%  "
# Row 10
$ws.Range("E10").Value = "'  -0.97%  "
# Row 11
$ws.Range("D11").Value = "'0.07790"
$ws.Range("E11").Value = "'  +1.20%  "
# Row 12
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.728.20"
$ws.Range("E12").Value = "'  +4.81%  "
# Row 13
$ws.Range("B13").Value = "'Polkadot"
$ws.Range("C13").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.503"
$ws.Range("E13").Value = "'  +2.51%  "
# Row 14
$ws.Range("D14").Value = "'0.5481"
$ws.Range("E14").Value = "'  +1.36%  "
# Row 15
$ws.Range("D15").Value = "'0.0₅8227"
$ws.Range("E15").Value = "'  +2.12%  "
# Row 16
$ws.Range("D16").Value = "'65.37"
# Row 17
$ws.Range("D17").Value = "'26.131.50"
$ws.Range("E17").Value = "'  -0.23%  "
# Row 18
$ws.Range("E18").Value = "'  -0.34%  "
# Row 19
$ws.Range("D19").Value = "'4.578"
$ws.Range("E19").Value = "'  -0.59%  "
# Row 20
$ws.Range("D20").Value = "'191.71"
$ws.Range("E20").Value = "'  +0.52%  "
# Row 21
$ws.Range("D21").Value = "'10.07"
$ws.Range("E21").Value = "'  +0.73%  "
# Row 22
$ws.Range("D22").Value = "'6.044"
$ws.Range("E22").Value = "'  +0.23%  "
# Row 23
$ws.Range("E23").Value = "'  -0.42%  "
# Row 24
$ws.Range("D24").Value = "'142.11"
$ws.Range("E24").Value = "'  +1.82%  "
# Row 25
$ws.Range("E25").Value = "'  +2.59%  "
# Row 26
$ws.Range("D26").Value = "'7.268"
$ws.Range("E26").Value = "'  +2.05%  "
# Row 27
$ws.Range("E27").Value = "'  +1.25%  "
# Row 28
$ws.Range("E28").Value = "'  +1.69%  "
# Row 29
$ws.Range("D29").Value = "'0.05914"
$ws.Range("E29").Value = "'  -0.85%  "
# Row 30
$ws.Range("D30").Value = "'1.278"
$ws.Range("E30").Value = "'  +0.48%  "
# Row 31
$ws.Range("D31").Value = "'3.526"
$ws.Range("E31").Value = "'  -1.19%  "
# Row 32
$ws.Range("D32").Value = "'3.253"
$ws.Range("E32").Value = "'  +0.44%  "
# Row 33
$ws.Range("E33").Value = "'  -1.49%  "
# Row 34
$ws.Range("D34").Value = "'0.9536"
$ws.Range("E34").Value = "'  -0.32%  "
# Row 35
$ws.Range("D35").Value = "'2.786"
$ws.Range("E35").Value = "'  +0.31%  "
# Row 36
$ws.Range("D36").Value = "'2.412"
$ws.Range("E36").Value = "'  -0.36%  "
# Row 37
$ws.Range("D37").Value = "'0.5709"
$ws.Range("E37").Value = "'  +1.34%  "
# Row 38
$ws.Range("D38").Value = "'0.01621"
$ws.Range("E38").Value = "'  +2.20%  "
# Row 39
$ws.Range("B39").Value = "'FraxShare"
$ws.Range("C39").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'5.783"
$ws.Range("E39").Value = "'  -3.58%  "
# Row 40
$ws.Range("B40").Value = "'TrustWalletToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.8484"
$ws.Range("E40").Value = "'  -0.85%  "
# Row 42
$ws.Range("D42").Value = "'1.029.18"
$ws.Range("E42").Value = "'  +1.29%  "
# Row 43
$ws.Range("D43").Value = "'102.85"
$ws.Range("E43").Value = "'  +2.99%  "
# Row 44
$ws.Range("D44").Value = "'1.801.59"
$ws.Range("E44").Value = "'  +0.12%  "
# Row 45
$ws.Range("D45").Value = "'57.31"
$ws.Range("E45").Value = "'  +1.31%  "
# Row 46
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "'  -0.86%  "
# Row 47
$ws.Range("E47").Value = "'  +2.12%  "
# Row 48
$ws.Range("E48").Value = "'  +1.62%  "
# Row 49
$ws.Range("E49").Value = "'  -0.29%  "
# Row 50
$ws.Range("D50").Value = "'7.855"
$ws.Range("E50").Value = "'  -0.96%  "
# Row 51
$ws.Range("D51").Value = "'0.09726"
$ws.Range("E51").Value = "'  +0.90%  "
